$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F4 831 -> 833, F5 916 -> 938
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 833
$ws1.Range("F5").Value = 938

# Sheet "全部类型" (sheet4): F5 831 -> 833, F6 916 -> 938
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 833
$ws4.Range("F6").Value = 938
